$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3167.5386
$ws.Range("J17").Value = 3167.5386
$ws.Range("L17").Value = 9502.6158
$ws.Range("N17").Value = -9838.6158
$ws.Range("H38").Value = 5032.5527
$ws.Range("I38").Value = 211
$ws.Range("J38").Value = 5936.5938
$ws.Range("K38").Value = 633
$ws.Range("L38").Value = 17809.7814
$ws.Range("M38").Value = -261
$ws.Range("N38").Value = -18553.7814
$ws.Range("H40").Value = 1663.5714
$ws.Range("I40").Value = 1860.4
$ws.Range("J40").Value = 1554.2222
$ws.Range("K40").Value = 1860.4
$ws.Range("L40").Value = 1554.2222
$ws.Range("M40").Value = -1685.4
$ws.Range("N40").Value = -1904.2222
$ws.Range("H43").Value = 5061159.5
$ws.Range("J43").Value = 9261526
$ws.Range("L43").Value = 9261526
$ws.Range("N43").Value = -9261664
$ws.Range("H58").Value = 928.1667
$ws.Range("J58").Value = 1950
$ws.Range("L58").Value = 5850
$ws.Range("N58").Value = -6150
$ws.Range("H62").Value = 12349454
$ws.Range("J62").Value = 4999.5
$ws.Range("L62").Value = 4999.5
$ws.Range("N62").Value = -6247.5
$ws.Range("H65").Value = 12349454
$ws.Range("J65").Value = 4999.5
$ws.Range("L65").Value = 24997.5
$ws.Range("N65").Value = -31237.5
$ws.Range("H86").Value = 14312.5
$ws.Range("I86").Value = 33966.668
$ws.Range("J86").Value = 2520
$ws.Range("K86").Value = 33966.668
$ws.Range("L86").Value = 2520
$ws.Range("M86").Value = -32843.668
$ws.Range("N86").Value = -4766
$ws.Range("H89").Value = 14312.5
$ws.Range("I89").Value = 33966.668
$ws.Range("J89").Value = 2520
$ws.Range("K89").Value = 169833.34
$ws.Range("L89").Value = 12600
$ws.Range("M89").Value = -164217.34
$ws.Range("N89").Value = -23832
$ws.Range("H138").Value = 1717.1
$ws.Range("I138").Value = 873.58826
$ws.Range("J138").Value = 1889.8674
$ws.Range("K138").Value = 2620.76478
$ws.Range("L138").Value = 5669.6022
$ws.Range("M138").Value = 2519.23522
$ws.Range("N138").Value = -15949.6022
$ws.Range("H140").Value = 33520
$ws.Range("J140").Value = 33520
$ws.Range("L140").Value = 33520
$ws.Range("N140").Value = -43880
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5433.1816
$ws.Range("I31").Value = 5433.1816
$ws.Range("K31").Value = 5433.1816
$ws.Range("M31").Value = -5139.1816
$ws.Range("H102").Value = 12821808
$ws.Range("I102").Value = 16667850
$ws.Range("J102").Value = 1670.3334
$ws.Range("K102").Value = 16667850
$ws.Range("L102").Value = 1670.3334
$ws.Range("M102").Value = -16666228
$ws.Range("N102").Value = -4914.3334
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H114").Value = 27300
$ws.Range("J114").Value = 27300
$ws.Range("L114").Value = 27300
$ws.Range("N114").Value = -35978
$ws.Range("H122").Value = 2223.2307
$ws.Range("I122").Value = 2191.389
$ws.Range("J122").Value = 2294.875
$ws.Range("K122").Value = 6574.167
$ws.Range("L122").Value = 6884.625
$ws.Range("M122").Value = -4124.167
$ws.Range("N122").Value = -11784.625
$ws.Range("H132").Value = 5060.364
$ws.Range("I132").Value = 7156
$ws.Range("J132").Value = 3862.8572
$ws.Range("K132").Value = 21468
$ws.Range("L132").Value = 11588.5716
$ws.Range("M132").Value = -18938
$ws.Range("N132").Value = -16648.5716
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 175400.75
$ws.Range("J22").Value = 233800.67
$ws.Range("L22").Value = 233800.67
$ws.Range("N22").Value = -234500.67
$ws.Range("H39").Value = 3200
$ws.Range("J39").Value = 3200
$ws.Range("L39").Value = 3200
$ws.Range("N39").Value = -3982
$ws.Range("H49").Value = 3200
$ws.Range("J49").Value = 3200
$ws.Range("L49").Value = 3200
$ws.Range("N49").Value = -3564
$ws.Range("H50").Value = 27000
$ws.Range("J50").Value = 27000
$ws.Range("L50").Value = 27000
$ws.Range("N50").Value = -28250
$ws.Range("H105").Value = 703.4167
$ws.Range("I105").Value = 672.4761999999999
$ws.Range("J105").Value = 920
$ws.Range("K105").Value = 672.4761999999999
$ws.Range("L105").Value = 920
$ws.Range("M105").Value = 1074.5238
$ws.Range("N105").Value = -4414
$ws.Range("H122").Value = 8437.429
$ws.Range("I122").Value = 8932.615
$ws.Range("K122").Value = 26797.845
$ws.Range("M122").Value = -24347.845
$ws.Range("H132").Value = 3433.5454
$ws.Range("I132").Value = 3474.2
$ws.Range("K132").Value = 10422.6
$ws.Range("M132").Value = -7892.599999999999
$ws.Range("H141").Value = 1510062
$ws.Range("J141").Value = 1510062
$ws.Range("L141").Value = 1510062
$ws.Range("N141").Value = -1520422
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 990.8182
$ws.Range("I2").Value = 48.42857
$ws.Range("J2").Value = 2640
$ws.Range("K2").Value = 290.57142
$ws.Range("L2").Value = 15840
$ws.Range("M2").Value = -177.57142
$ws.Range("N2").Value = -16066
$ws.Range("H7").Value = 448.64706
$ws.Range("I7").Value = 448.53333
$ws.Range("J7").Value = 449.5
$ws.Range("K7").Value = 1345.59999
$ws.Range("L7").Value = 1348.5
$ws.Range("M7").Value = -1233.59999
$ws.Range("N7").Value = -1572.5
$ws.Range("H117").Value = 763
$ws.Range("I117").Value = 372.2857
$ws.Range("K117").Value = 1116.8571
$ws.Range("M117").Value = 2325.1429
$ws.Range("H119").Value = 7015.8
$ws.Range("I119").Value = 2386
$ws.Range("K119").Value = 7158
$ws.Range("M119").Value = -2320
$ws.Range("H129").Value = 15433084
$ws.Range("I129").Value = 47619816
$ws.Range("J129").Value = 4167728.2
$ws.Range("K129").Value = 142859448
$ws.Range("L129").Value = 12503184.6
$ws.Range("M129").Value = -142854448
$ws.Range("N129").Value = -12513184.6
$ws.Range("H130").Value = 2328.25
$ws.Range("I130").Value = 1015
$ws.Range("J130").Value = 2766
$ws.Range("K130").Value = 3045
$ws.Range("L130").Value = 8298
$ws.Range("M130").Value = 1975
$ws.Range("N130").Value = -18338
$ws.Range("H131").Value = 10753750
$ws.Range("I131").Value = 142857540
$ws.Range("J131").Value = 1116.5814
$ws.Range("K131").Value = 428572620
$ws.Range("L131").Value = 3349.7442
$ws.Range("M131").Value = -428567580
$ws.Range("N131").Value = -13429.7442
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5609
$ws.Range("I80").Value = 4883.3335
$ws.Range("J80").Value = 6479.8
$ws.Range("K80").Value = 4883.3335
$ws.Range("L80").Value = 6479.8
$ws.Range("M80").Value = -3885.3335
$ws.Range("N80").Value = -8475.799999999999
$ws.Range("H83").Value = 5609
$ws.Range("I83").Value = 4883.3335
$ws.Range("J83").Value = 6479.8
$ws.Range("K83").Value = 24416.6675
$ws.Range("L83").Value = 32399
$ws.Range("M83").Value = -19424.6675
$ws.Range("N83").Value = -42383
$ws.Range("H122").Value = 2184.862
$ws.Range("I122").Value = 2332.9092
$ws.Range("J122").Value = 1719.5714
$ws.Range("K122").Value = 6998.7276
$ws.Range("L122").Value = 5158.7142
$ws.Range("M122").Value = -4548.7276
$ws.Range("N122").Value = -10058.7142
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 2000000
$ws.Range("I25").Value = 2000000
$ws.Range("K25").Value = 2000000
$ws.Range("M25").Value = -1999770
$ws.Range("H132").Value = 3893.0667
$ws.Range("I132").Value = 5199.6
$ws.Range("K132").Value = 15598.8
$ws.Range("M132").Value = -13068.8
$ws.Range("H136").Value = 1087.5
$ws.Range("I136").Value = 1157.1428
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 3471.4284
$ws.Range("L136").Value = 1800
$ws.Range("M136").Value = -921.4284000000002
$ws.Range("N136").Value = -6900
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2540.5293
$ws.Range("I132").Value = 2085.0715
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 6255.2145
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -3725.2145
$ws.Range("N132").Value = -19058
$ws.Range("H136").Value = 1024.8684
$ws.Range("I136").Value = 932.7931
$ws.Range("J136").Value = 1321.5555
$ws.Range("K136").Value = 2798.3793
$ws.Range("L136").Value = 3964.6665
$ws.Range("M136").Value = -248.3793000000001
$ws.Range("N136").Value = -9064.666499999999
